$d = $word.ActiveDocument

# --- 1. Merge the two runs of paragraph 2 ("PARAGRAFO " + "2") into a single run ---
# Scope the Find to paragraph 2 only, and replace the text with itself so the
# engine re-serializes it as one run (no other paragraph's "PARAGRAFO" text is touched).
$p2 = $d.Paragraphs.Item(2).Range
$p2.Find.Execute("PARAGRAFO ", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "PARAGRAFO ", 2) | Out-Null

# --- 2. Add "Ramo1: " before the _GoBack bookmark and "inizio modifiche ramo1" after it ---
# Paragraph 3 is the paragraph that only holds the bookmarkStart/bookmarkEnd pair.
$p3 = $d.Paragraphs.Item(3).Range
$bmPos = $p3.Start

# Insert the "after bookmark" text first (while the bookmark still sits at $bmPos),
# using InsertAfter on a collapsed range so the new run lands after the bookmark tags.
$afterAnchor = $d.Range($bmPos, $bmPos)
$afterAnchor.InsertAfter("inizio modifiche ramo1") | Out-Null

# Now insert the "before bookmark" text at the same original position; InsertBefore
# places the new run ahead of the (still unmoved) bookmark tags.
$beforeAnchor = $d.Range($bmPos, $bmPos)
$beforeAnchor.InsertBefore("Ramo1: ") | Out-Null

# --- 3. Append two additional empty paragraphs before the existing trailing empty paragraph ---
# InsertParagraphBefore() leaves a stray empty run behind in the new paragraph, so we
# type a placeholder character into it and delete that character again; this makes the
# run disappear entirely, leaving a truly empty <w:p/>.
for ($n = 0; $n -lt 2; $n++) {
    $lastIdx = $d.Paragraphs.Count
    $origLast = $d.Paragraphs.Item($lastIdx)
    $origLast.Range.InsertParagraphBefore() | Out-Null

    $newP = $d.Paragraphs.Item($lastIdx)
    $s = $newP.Range.Start
    $placeholder = $d.Range($s, $s)
    $placeholder.InsertBefore("X") | Out-Null

    $newP2 = $d.Paragraphs.Item($lastIdx)
    $cleanup = $d.Range($newP2.Range.Start, $newP2.Range.Start + 1)
    $cleanup.Delete() | Out-Null
}

Write-Output "ok"
